$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 89049.42999999999
$ws.Range("I86").Value = 176977.86
$ws.Range("J86").Value = 1121
$ws.Range("K86").Value = 176977.86
$ws.Range("L86").Value = 1121
$ws.Range("M86").Value = -175854.86
$ws.Range("N86").Value = -3367
$ws.Range("H89").Value = 89049.42999999999
$ws.Range("I89").Value = 176977.86
$ws.Range("J89").Value = 1121
$ws.Range("K89").Value = 884889.2999999999
$ws.Range("L89").Value = 5605
$ws.Range("M89").Value = -879273.2999999999
$ws.Range("N89").Value = -16837
$ws.Range("H112").Value = 1360.75
$ws.Range("J112").Value = 1370.4073
$ws.Range("L112").Value = 4111.2219
$ws.Range("N112").Value = -6327.2219
$ws.Range("H132").Value = 6804161
$ws.Range("I132").Value = 7937862
$ws.Range("K132").Value = 23813586
$ws.Range("M132").Value = -23811056
$ws.Range("H137").Value = 1142.1628
$ws.Range("I137").Value = 916.3200000000001
$ws.Range("K137").Value = 2748.96
$ws.Range("M137").Value = -198.96
$ws.Range("H138").Value = 1895.7941
$ws.Range("I138").Value = 1640.4694
$ws.Range("J138").Value = 2554.2632
$ws.Range("K138").Value = 4921.4082
$ws.Range("L138").Value = 7662.7896
$ws.Range("M138").Value = 218.5918000000001
$ws.Range("N138").Value = -17942.7896

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 12499.5
$ws.Range("I22").Value = 4999
$ws.Range("J22").Value = 20000
$ws.Range("K22").Value = 4999
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = -4700
$ws.Range("N22").Value = -20598
$ws.Range("H32").Value = 2892.5276
$ws.Range("I32").Value = 2467.7908
$ws.Range("K32").Value = 2467.7908
$ws.Range("M32").Value = -2180.7908
$ws.Range("H45").Value = 1711.3043
$ws.Range("I45").Value = 1603.0714
$ws.Range("J45").Value = 1879.6666
$ws.Range("K45").Value = 1603.0714
$ws.Range("L45").Value = 1879.6666
$ws.Range("M45").Value = -1226.0714
$ws.Range("N45").Value = -2633.6666
$ws.Range("H61").Value = 43481016
$ws.Range("I61").Value = 25002068
$ws.Range("K61").Value = 25002068
$ws.Range("M61").Value = -25001856
$ws.Range("H122").Value = 1404.625
$ws.Range("I122").Value = 1282.4375
$ws.Range("K122").Value = 3847.3125
$ws.Range("M122").Value = -1397.3125
$ws.Range("H136").Value = 43481016
$ws.Range("I136").Value = 25002068
$ws.Range("K136").Value = 75006204
$ws.Range("M136").Value = -75003654

# --- BSM sheet updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 575089.6
$ws.Range("I86").Value = 837740.2
$ws.Range("J86").Value = 224888.89
$ws.Range("K86").Value = 837740.2
$ws.Range("L86").Value = 224888.89
$ws.Range("M86").Value = -836617.2
$ws.Range("N86").Value = -227134.89
$ws.Range("H89").Value = 575089.6
$ws.Range("I89").Value = 837740.2
$ws.Range("J89").Value = 224888.89
$ws.Range("K89").Value = 4188701
$ws.Range("L89").Value = 1124444.45
$ws.Range("M89").Value = -4183085
$ws.Range("N89").Value = -1135676.45
$ws.Range("H134").Value = 10361.071
$ws.Range("I134").Value = 11004.583
$ws.Range("K134").Value = 33013.749
$ws.Range("M134").Value = -30478.749

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2306021.5
$ws.Range("I31").Value = 3969692.2
$ws.Range("J31").Value = 2477.3076
$ws.Range("K31").Value = 3969692.2
$ws.Range("L31").Value = 2477.3076
$ws.Range("M31").Value = -3969397.2
$ws.Range("N31").Value = -3067.3076
$ws.Range("H34").Value = 2306021.5
$ws.Range("I34").Value = 3969692.2
$ws.Range("J34").Value = 2477.3076
$ws.Range("K34").Value = 3969692.2
$ws.Range("L34").Value = 2477.3076
$ws.Range("M34").Value = -3969490.2
$ws.Range("N34").Value = -2881.3076
$ws.Range("H58").Value = 2290054
$ws.Range("I58").Value = 3624393
$ws.Range("J58").Value = 2615.7144
$ws.Range("K58").Value = 3624393
$ws.Range("L58").Value = 2615.7144
$ws.Range("M58").Value = -3624190
$ws.Range("N58").Value = -3021.7144
$ws.Range("H132").Value = 1582.9565
$ws.Range("I132").Value = 1193.6216
$ws.Range("J132").Value = 3183.5557
$ws.Range("K132").Value = 3580.8648
$ws.Range("L132").Value = 9550.667099999999
$ws.Range("M132").Value = -1050.8648
$ws.Range("N132").Value = -14610.6671
$ws.Range("H136").Value = 2290054
$ws.Range("I136").Value = 3624393
$ws.Range("J136").Value = 2615.7144
$ws.Range("K136").Value = 10873179
$ws.Range("L136").Value = 7847.1432
$ws.Range("M136").Value = -10870629
$ws.Range("N136").Value = -12947.1432

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 797
$ws.Range("J117").Value = 899.8
$ws.Range("L117").Value = 2699.4
$ws.Range("N117").Value = -9583.4
$ws.Range("H140").Value = 2075.0186
$ws.Range("I140").Value = 1026
$ws.Range("K140").Value = 3078
$ws.Range("M140").Value = 2102

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2477.9092
$ws.Range("I102").Value = 2508.8
$ws.Range("J102").Value = 2411.7144
$ws.Range("K102").Value = 2508.8
$ws.Range("L102").Value = 2411.7144
$ws.Range("M102").Value = -886.8000000000002
$ws.Range("N102").Value = -5655.7144
$ws.Range("H122").Value = 1300.3572
$ws.Range("I122").Value = 1224.8572
$ws.Range("J122").Value = 1526.8572
$ws.Range("K122").Value = 3674.5716
$ws.Range("L122").Value = 4580.571599999999
$ws.Range("M122").Value = -1224.5716
$ws.Range("N122").Value = -9480.571599999999
$ws.Range("H132").Value = 688363.7
$ws.Range("I132").Value = 855875
$ws.Range("J132").Value = 3090.0908
$ws.Range("K132").Value = 2567625
$ws.Range("L132").Value = 9270.2724
$ws.Range("M132").Value = -2565095
$ws.Range("N132").Value = -14330.2724
$ws.Range("H138").Value = 50268.8
$ws.Range("J138").Value = 50268.8
$ws.Range("L138").Value = 50268.8
$ws.Range("N138").Value = -60548.8

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1825.2858
$ws.Range("I46").Value = 1060.25
$ws.Range("K46").Value = 1060.25
$ws.Range("M46").Value = -872.25
$ws.Range("H132").Value = 2149.1428
$ws.Range("I132").Value = 1566.76
$ws.Range("J132").Value = 3605.1
$ws.Range("K132").Value = 4700.28
$ws.Range("L132").Value = 10815.3
$ws.Range("M132").Value = -2170.28
$ws.Range("N132").Value = -15875.3
$ws.Range("H136").Value = 2986.2432
$ws.Range("I136").Value = 1956.4
$ws.Range("K136").Value = 5869.200000000001
$ws.Range("M136").Value = -3319.200000000001

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 15000
$ws.Range("J22").Value = 15000
$ws.Range("L22").Value = 15000
$ws.Range("N22").Value = -15586
$ws.Range("H46").Value = 43975.668
$ws.Range("J46").Value = 43975.668
$ws.Range("L46").Value = 43975.668
$ws.Range("N46").Value = -44437.668
$ws.Range("H126").Value = 3463.1875
$ws.Range("I126").Value = 2642
$ws.Range("J126").Value = 3955.9
$ws.Range("K126").Value = 7926
$ws.Range("L126").Value = 11867.7
$ws.Range("M126").Value = -5456
$ws.Range("N126").Value = -16807.7
$ws.Range("H132").Value = 1443.8837
$ws.Range("I132").Value = 1175.4073
$ws.Range("J132").Value = 1896.9375
$ws.Range("K132").Value = 3526.2219
$ws.Range("L132").Value = 5690.8125
$ws.Range("M132").Value = -996.2219000000005
$ws.Range("N132").Value = -10750.8125
$ws.Range("H134").Value = 43975.668
$ws.Range("J134").Value = 43975.668
$ws.Range("L134").Value = 131927.004
$ws.Range("N134").Value = -136997.004
$ws.Range("H136").Value = 12079506
$ws.Range("I136").Value = 14247160
$ws.Range("K136").Value = 42741480
$ws.Range("M136").Value = -42738930
$ws.Range("H139").Value = 69667.71000000001
$ws.Range("J139").Value = 69667.71000000001
$ws.Range("L139").Value = 69667.71000000001
$ws.Range("N139").Value = -79947.71000000001
